$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update numeric values in column C
$ws.Range("C3").Value = 28
$ws.Range("C4").Value = 26
$ws.Range("C5").Value = 33
$ws.Range("C6").Value = 22
$ws.Range("C7").Value = 16
$ws.Range("C8").Value = 21
$ws.Range("C9").Value = 33
$ws.Range("C10").Value = 30
$ws.Range("C11").Value = 27
$ws.Range("C12").Value = 38
$ws.Range("C13").Value = 32
$ws.Range("C14").Value = 24
$ws.Range("C15").Value = 28
$ws.Range("C16").Value = 22
$ws.Range("C17").Value = 26
$ws.Range("C18").Value = 24

# Update text values in column B
$ws.Range("B9").Value = "<it>"
$ws.Range("B13").Value = "<victer>"
$ws.Range("B14").Value = "<alt>"
$ws.Range("B16").Value = "<yankee>"
$ws.Range("B17").Value = "<sen>"
$ws.Range("B18").Value = "<ha>"
